# Board.xlsx update - "divisão das tarefas relativas ao relatório"
# Fills in the task/owner columns (B7:C15) for the "Relatório" section rows
# that were previously left blank, and updates the current sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B first: the task names for the "Relatório" breakdown (rows 7-15).
$ws.Range("B7").Value = "Abstract"
$ws.Range("B8").Value = "Introdução"
$ws.Range("B9").Value = "State-of-the-art"
$ws.Range("B10").Value = "Problema"
$ws.Range("B11").Value = "Solução"
$ws.Range("B12").Value = "Conclusão"
$ws.Range("B13").Value = "Referencias"
$ws.Range("B14").Value = "Revisão Ortográfica"
$ws.Range("B15").Value = "Revisão da Organização"

# Then column C: who each task is assigned to.
$ws.Range("C7").Value = "Ricardo Osório"
$ws.Range("C8").Value = "João Cardoso"
$ws.Range("C9").Value = "Ana Leite"
$ws.Range("C10").Value = "Ana Leite"
$ws.Range("C11").Value = "Guilherme Ferreira"
$ws.Range("C12").Value = "Ricardo Catalão"
$ws.Range("C13").Value = "Todos"
$ws.Range("C14").Value = "João Cardoso e Ricardo Catalão"
$ws.Range("C15").Value = "Ana Leite e Guilherme Ferreira"

# Leave the sheet scrolled/selected where the author ended up editing.
$ws.Range("C16").Select()
